$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 889.125
$ws.Range("I6").Value = 710.8333
$ws.Range("J6").Value = 1424
$ws.Range("K6").Value = 2132.4999
$ws.Range("L6").Value = 4272
$ws.Range("M6").Value = -2020.4999
$ws.Range("N6").Value = -4496
$ws.Range("H74").Value = 8266.666999999999
$ws.Range("I74").Value = 8266.666999999999
$ws.Range("K74").Value = 8266.666999999999
$ws.Range("M74").Value = -7330.666999999999
$ws.Range("H76").Value = 5449.875
$ws.Range("I76").Value = 4599.6
$ws.Range("J76").Value = 6867
$ws.Range("K76").Value = 4599.6
$ws.Range("L76").Value = 6867
$ws.Range("M76").Value = -4284.6
$ws.Range("N76").Value = -7497
$ws.Range("H77").Value = 8266.666999999999
$ws.Range("I77").Value = 8266.666999999999
$ws.Range("K77").Value = 41333.335
$ws.Range("M77").Value = -36653.335
$ws.Range("H79").Value = 5449.875
$ws.Range("I79").Value = 4599.6
$ws.Range("J79").Value = 6867
$ws.Range("K79").Value = 4599.6
$ws.Range("L79").Value = 6867
$ws.Range("M79").Value = -3507.6
$ws.Range("N79").Value = -9051
$ws.Range("H132").Value = 2321.5117
$ws.Range("I132").Value = 2475.0527
$ws.Range("K132").Value = 7425.158100000001
$ws.Range("M132").Value = -4895.158100000001
$ws.Range("H138").Value = 2789.2632
$ws.Range("I138").Value = 1153.9375
$ws.Range("J138").Value = 3978.5908
$ws.Range("K138").Value = 3461.8125
$ws.Range("L138").Value = 11935.7724
$ws.Range("M138").Value = 1678.1875
$ws.Range("N138").Value = -22215.7724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4298.3486
$ws.Range("I61").Value = 2848
$ws.Range("K61").Value = 2848
$ws.Range("M61").Value = -2636
$ws.Range("H63").Value = 2790
$ws.Range("I63").Value = 2790
$ws.Range("K63").Value = 2790
$ws.Range("M63").Value = -2104
$ws.Range("H66").Value = 2790
$ws.Range("I66").Value = 2790
$ws.Range("K66").Value = 13950
$ws.Range("M66").Value = -10518
$ws.Range("H88").Value = 2874.75
$ws.Range("I88").Value = 2999
$ws.Range("J88").Value = 2833.3333
$ws.Range("K88").Value = 2999
$ws.Range("L88").Value = 2833.3333
$ws.Range("M88").Value = -2593
$ws.Range("N88").Value = -3645.3333
$ws.Range("H91").Value = 2874.75
$ws.Range("I91").Value = 2999
$ws.Range("J91").Value = 2833.3333
$ws.Range("K91").Value = 2999
$ws.Range("L91").Value = 2833.3333
$ws.Range("M91").Value = -1595
$ws.Range("N91").Value = -5641.3333
$ws.Range("H97").Value = 985.7857
$ws.Range("I97").Value = 612.8182
$ws.Range("J97").Value = 2353.3333
$ws.Range("K97").Value = 612.8182
$ws.Range("L97").Value = 2353.3333
$ws.Range("M97").Value = -116.8182
$ws.Range("N97").Value = -3345.3333
$ws.Range("H132").Value = 2308.516
$ws.Range("I132").Value = 2218.4614
$ws.Range("J132").Value = 2461.2173
$ws.Range("K132").Value = 6655.3842
$ws.Range("L132").Value = 7383.651899999999
$ws.Range("M132").Value = -4125.3842
$ws.Range("N132").Value = -12443.6519
$ws.Range("H136").Value = 4298.3486
$ws.Range("I136").Value = 2848
$ws.Range("K136").Value = 8544
$ws.Range("M136").Value = -5994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1613.591
$ws.Range("I20").Value = 1769.875
$ws.Range("J20").Value = 1196.8334
$ws.Range("K20").Value = 1769.875
$ws.Range("L20").Value = 1196.8334
$ws.Range("M20").Value = -1522.875
$ws.Range("N20").Value = -1690.8334
$ws.Range("H86").Value = 35780460
$ws.Range("I86").Value = 58878150
$ws.Range("J86").Value = 84023.91
$ws.Range("K86").Value = 58878150
$ws.Range("L86").Value = 84023.91
$ws.Range("M86").Value = -58877027
$ws.Range("N86").Value = -86269.91
$ws.Range("H89").Value = 35780460
$ws.Range("I89").Value = 58878150
$ws.Range("J89").Value = 84023.91
$ws.Range("K89").Value = 294390750
$ws.Range("L89").Value = 420119.55
$ws.Range("M89").Value = -294385134
$ws.Range("N89").Value = -431351.55
$ws.Range("H105").Value = 62503124
$ws.Range("I105").Value = 71431570
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 71431570
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -71429823
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2173.98
$ws.Range("I31").Value = 1749.2307
$ws.Range("J31").Value = 2634.125
$ws.Range("K31").Value = 1749.2307
$ws.Range("L31").Value = 2634.125
$ws.Range("M31").Value = -1454.2307
$ws.Range("N31").Value = -3224.125
$ws.Range("H34").Value = 2173.98
$ws.Range("I34").Value = 1749.2307
$ws.Range("J34").Value = 2634.125
$ws.Range("K34").Value = 1749.2307
$ws.Range("L34").Value = 2634.125
$ws.Range("M34").Value = -1547.2307
$ws.Range("N34").Value = -3038.125
$ws.Range("H58").Value = 2898.5925
$ws.Range("I58").Value = 2428.3809
$ws.Range("J58").Value = 4544.3335
$ws.Range("K58").Value = 2428.3809
$ws.Range("L58").Value = 4544.3335
$ws.Range("M58").Value = -2225.3809
$ws.Range("N58").Value = -4950.3335
$ws.Range("H62").Value = 84021.16
$ws.Range("I62").Value = 128685.625
$ws.Range("J62").Value = 12558
$ws.Range("K62").Value = 128685.625
$ws.Range("L62").Value = 12558
$ws.Range("M62").Value = -128061.625
$ws.Range("N62").Value = -13806
$ws.Range("H65").Value = 84021.16
$ws.Range("I65").Value = 128685.625
$ws.Range("J65").Value = 12558
$ws.Range("K65").Value = 643428.125
$ws.Range("L65").Value = 62790
$ws.Range("M65").Value = -640308.125
$ws.Range("N65").Value = -69030
$ws.Range("H136").Value = 2898.5925
$ws.Range("I136").Value = 2428.3809
$ws.Range("J136").Value = 4544.3335
$ws.Range("K136").Value = 7285.1427
$ws.Range("L136").Value = 13633.0005
$ws.Range("M136").Value = -4735.1427
$ws.Range("N136").Value = -18733.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 330
$ws.Range("I86").Value = 300
$ws.Range("J86").Value = 350
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 1050
$ws.Range("M86").Value = 286
$ws.Range("N86").Value = -3422
$ws.Range("H89").Value = 330
$ws.Range("I89").Value = 300
$ws.Range("J89").Value = 350
$ws.Range("K89").Value = 2700
$ws.Range("L89").Value = 3150
$ws.Range("M89").Value = 3228
$ws.Range("N89").Value = -15006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 6963.7144
$ws.Range("I29").Value = 4819.8
$ws.Range("K29").Value = 4819.8
$ws.Range("M29").Value = -4529.8
$ws.Range("H80").Value = 3470.182
$ws.Range("I80").Value = 3527.2
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 3527.2
$ws.Range("L80").Value = 2900
$ws.Range("M80").Value = -2529.2
$ws.Range("N80").Value = -4896
$ws.Range("H83").Value = 3470.182
$ws.Range("I83").Value = 3527.2
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 17636
$ws.Range("L83").Value = 14500
$ws.Range("M83").Value = -12644
$ws.Range("N83").Value = -24484
$ws.Range("H97").Value = 25665982
$ws.Range("I97").Value = 29439876
$ws.Range("J97").Value = 3501.8
$ws.Range("K97").Value = 29439876
$ws.Range("L97").Value = 3501.8
$ws.Range("M97").Value = -29439380
$ws.Range("N97").Value = -4493.8
$ws.Range("H102").Value = 31536.895
$ws.Range("I102").Value = 37802.31
$ws.Range("J102").Value = 11348.333
$ws.Range("K102").Value = 37802.31
$ws.Range("L102").Value = 11348.333
$ws.Range("M102").Value = -36180.31
$ws.Range("N102").Value = -14592.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 2197
$ws.Range("I23").Value = 2197
$ws.Range("K23").Value = 2197
$ws.Range("M23").Value = -1967
$ws.Range("H40").Value = 5028
$ws.Range("I40").Value = 3281.04
$ws.Range("J40").Value = 6534
$ws.Range("K40").Value = 3281.04
$ws.Range("L40").Value = 6534
$ws.Range("M40").Value = -3145.04
$ws.Range("N40").Value = -6806
$ws.Range("H46").Value = 3060.9312
$ws.Range("I46").Value = 1175.9231
$ws.Range("J46").Value = 4592.5
$ws.Range("K46").Value = 1175.9231
$ws.Range("L46").Value = 4592.5
$ws.Range("M46").Value = -987.9231
$ws.Range("N46").Value = -4968.5
$ws.Range("H55").Value = 1344.7858
$ws.Range("I55").Value = 143.16667
$ws.Range("J55").Value = 3507.7
$ws.Range("K55").Value = 143.16667
$ws.Range("L55").Value = 3507.7
$ws.Range("M55").Value = 29.83332999999999
$ws.Range("N55").Value = -3853.7
$ws.Range("H93").Value = 1321.5
$ws.Range("I93").Value = 1153.1428
$ws.Range("K93").Value = 1153.1428
$ws.Range("M93").Value = 94.85719999999992
$ws.Range("H132").Value = 4824.079
$ws.Range("J132").Value = 6967.625
$ws.Range("L132").Value = 20902.875
$ws.Range("N132").Value = -25962.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 13454.546
$ws.Range("I30").Value = 8615.385
$ws.Range("J30").Value = 31428.572
$ws.Range("K30").Value = 8615.385
$ws.Range("L30").Value = 31428.572
$ws.Range("M30").Value = -8508.385
$ws.Range("N30").Value = -31642.572
$ws.Range("H32").Value = 25000
$ws.Range("J32").Value = 25000
$ws.Range("L32").Value = 25000
$ws.Range("N32").Value = -25634
$ws.Range("H34").Value = 45000
$ws.Range("J34").Value = 45000
$ws.Range("L34").Value = 45000
$ws.Range("M34").Value = -45406
$ws.Range("H136").Value = 59224.65
$ws.Range("I136").Value = 44776.875
$ws.Range("J136").Value = 93899.3
$ws.Range("K136").Value = 134330.625
$ws.Range("L136").Value = 281697.9
$ws.Range("M136").Value = -131780.625
$ws.Range("N136").Value = -286797.9
